$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Price (D) and Volume(1h) (E) columns for rows 2-51 with freshly
# scraped crypto data. Values that look like plain numbers (e.g. "0.9994")
# are assigned with a leading apostrophe so Excel keeps them as literal text
# (matching the source file's inline-string cells) instead of reinterpreting
# them as numeric values; the cell style is then reset to "Normal" so no
# extra formatting (e.g. a Text number format) is left behind.

$ws.Range("D2").Value = '30.803.61'
$ws.Range("E2").Value = '  +0.98%  '
$ws.Range("D3").Value = '1.956.85'
$ws.Range("E3").Value = '  +3.79%  '
$ws.Range("D4").Value = "'0.9994"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").Value = "'251.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.54%  '
$ws.Range("D6").Value = "'0.5992"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +27.90%  '
$ws.Range("D7").Value = "'0.9989"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.22%  '
$ws.Range("D8").Value = "'0.3149"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +8.76%  '
$ws.Range("D9").Value = "'24.87"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +11.80%  '
$ws.Range("D10").Value = "'0.06866"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.82%  '
$ws.Range("D11").Value = "'0.8108"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +11.33%  '
$ws.Range("D12").Value = "'101.53"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +6.41%  '
$ws.Range("D13").Value = "'0.07974"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.96%  '
$ws.Range("D14").Value = '1.933.88'
$ws.Range("E14").Value = '  +2.57%  '
$ws.Range("D15").Value = "'5.354"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.41%  '
$ws.Range("D16").Value = "'280.47"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.48%  '
$ws.Range("D17").Value = '30.815.66'
$ws.Range("E17").Value = '  +1.04%  '
$ws.Range("D18").Value = "'13.76"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +5.75%  '
$ws.Range("D19").Value = "'0.000007700"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.93%  '
$ws.Range("D20").Value = "'5.606"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +6.54%  '
$ws.Range("D21").Value = '2.190.01'
$ws.Range("E21").Value = '  +2.88%  '
$ws.Range("D22").Value = "'1.0000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.07%  '
$ws.Range("D23").Value = "'0.9993"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.25%  '
$ws.Range("D24").Value = "'6.639"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +6.11%  '
$ws.Range("D25").Value = "'9.493"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.46%  '
$ws.Range("D26").Value = "'165.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.09%  '
$ws.Range("D27").Value = "'19.51"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.01%  '
$ws.Range("D28").Value = "'2.107"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +11.31%  '
$ws.Range("D29").Value = "'0.1132"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +16.58%  '
$ws.Range("D30").Value = "'1.359"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.91%  '
$ws.Range("D31").Value = "'1.555"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.81%  '
$ws.Range("D32").Value = "'4.461"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.23%  '
$ws.Range("D33").Value = "'4.359"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.50%  '
$ws.Range("D34").Value = "'0.05000"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.95%  '
$ws.Range("D35").Value = "'1.194"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.18%  '
$ws.Range("D36").Value = "'0.7203"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.87%  '
$ws.Range("D37").Value = "'2.727"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.26%  '
$ws.Range("D38").Value = "'0.01966"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.95%  '
$ws.Range("D39").Value = "'2.927"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.36%  '
$ws.Range("D40").Value = "'77.91"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.88%  '
$ws.Range("D41").Value = "'6.460"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.29%  '
$ws.Range("D42").Value = "'0.4523"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.36%  '
$ws.Range("D43").Value = "'2.021"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.77%  '
$ws.Range("D44").Value = "'0.8466"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.72%  '
$ws.Range("D45").Value = "'0.9991"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.17%  '
$ws.Range("D46").Value = "'10.07"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.85%  '
$ws.Range("D47").Value = "'102.58"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.22%  '
$ws.Range("D48").Value = "'7.308"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.89%  '
$ws.Range("D49").Value = "'36.10"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.61%  '
$ws.Range("D50").Value = "'0.4154"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.65%  '
$ws.Range("D51").Value = "'915.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.22%  '
